$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2858.4644
$ws.Range("J17").Value = 2419.6538
$ws.Range("L17").Value = 7258.9614
$ws.Range("N17").Value = -7594.9614
$ws.Range("H81").Value = 20149
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 20149
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H118").Value = 491.14285
$ws.Range("I118").Value = 491.14285
$ws.Range("K118").Value = 1473.42855
$ws.Range("M118").Value = 183.5714499999999
$ws.Range("H125").Value = 379
$ws.Range("I125").Value = 386.18182
$ws.Range("J125").Value = 300
$ws.Range("K125").Value = 3475.63638
$ws.Range("L125").Value = 2700
$ws.Range("M125").Value = -1015.63638
$ws.Range("N125").Value = -7620
$ws.Range("H140").Value = 64272.527
$ws.Range("J140").Value = 64272.527
$ws.Range("L140").Value = 64272.527
$ws.Range("N140").Value = -74632.527
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1915.7894
$ws.Range("I61").Value = 1854
$ws.Range("K61").Value = 1854
$ws.Range("M61").Value = -1642
$ws.Range("H74").Value = 1090.9117
$ws.Range("I74").Value = 603.38464
$ws.Range("K74").Value = 603.38464
$ws.Range("M74").Value = 270.61536
$ws.Range("H77").Value = 1090.9117
$ws.Range("I77").Value = 603.38464
$ws.Range("K77").Value = 3016.9232
$ws.Range("M77").Value = 1351.0768
$ws.Range("H132").Value = 1413.7273
$ws.Range("I132").Value = 1413.7273
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4241.1819
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1711.1819
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 1915.7894
$ws.Range("I136").Value = 1854
$ws.Range("K136").Value = 5562
$ws.Range("M136").Value = -3012
$ws.Range("H139").Value = 41593.332
$ws.Range("J139").Value = 41593.332
$ws.Range("L139").Value = 41593.332
$ws.Range("N139").Value = -51873.332
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 34950
$ws.Range("J103").Value = 34950
$ws.Range("L103").Value = 34950
$ws.Range("N103").Value = -37294
$ws.Range("H135").Value = 46853.168
$ws.Range("J135").Value = 46853.168
$ws.Range("L135").Value = 46853.168
$ws.Range("N135").Value = -56993.168
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 75801.25
$ws.Range("I16").Value = 86487.14
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 86487.14
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -86200.14
$ws.Range("N16").Value = -1574
$ws.Range("H31").Value = 2442.625
$ws.Range("I31").Value = 2255.4546
$ws.Range("K31").Value = 2255.4546
$ws.Range("M31").Value = -1960.4546
$ws.Range("H34").Value = 2442.625
$ws.Range("I34").Value = 2255.4546
$ws.Range("K34").Value = 2255.4546
$ws.Range("M34").Value = -2053.4546
$ws.Range("H43").Value = 14000
$ws.Range("J43").Value = 14000
$ws.Range("L43").Value = 14000
$ws.Range("N43").Value = -14368
$ws.Range("H87").Value = 34000
$ws.Range("J87").Value = 34000
$ws.Range("L87").Value = 34000
$ws.Range("N87").Value = -36372
$ws.Range("H90").Value = 34000
$ws.Range("J90").Value = 34000
$ws.Range("L90").Value = 102000
$ws.Range("N90").Value = -113856
$ws.Range("H101").Value = 14000
$ws.Range("J101").Value = 14000
$ws.Range("L101").Value = 14000
$ws.Range("N101").Value = -20490
$ws.Range("H113").Value = 75801.25
$ws.Range("I113").Value = 86487.14
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 86487.14
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = -84317.14
$ws.Range("N113").Value = -5340
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 14179.812
$ws.Range("J131").Value = 14432.692
$ws.Range("L131").Value = 43298.076
$ws.Range("N131").Value = -53378.076
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 20406000
$ws.Range("I10").Value = 20406000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 20406000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -20405831
$ws.Range("N10").ClearContents()
$ws.Range("H15").Value = 34998
$ws.Range("J15").Value = 34998
$ws.Range("L15").Value = 34998
$ws.Range("N15").Value = -35574
$ws.Range("H81").Value = 34998
$ws.Range("J81").Value = 34998
$ws.Range("L81").Value = 34998
$ws.Range("N81").Value = -36994
$ws.Range("H84").Value = 34998
$ws.Range("J84").Value = 34998
$ws.Range("L84").Value = 104994
$ws.Range("N84").Value = -114978
$ws.Range("H102").Value = 2283.8333
$ws.Range("I102").Value = 2132.1875
$ws.Range("J102").Value = 3497
$ws.Range("K102").Value = 2132.1875
$ws.Range("L102").Value = 3497
$ws.Range("M102").Value = -510.1875
$ws.Range("N102").Value = -6741
$ws.Range("H127").Value = 37217.5
$ws.Range("J127").Value = 37217.5
$ws.Range("L127").Value = 37217.5
$ws.Range("N127").Value = -47137.5
$ws.Range("H140").Value = 47533.223
$ws.Range("J140").Value = 47533.223
$ws.Range("L140").Value = 47533.223
$ws.Range("N140").Value = -57893.223
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10266.5
$ws.Range("I40").Value = 11718.454
$ws.Range("K40").Value = 11718.454
$ws.Range("M40").Value = -11582.454
$ws.Range("H93").Value = 15152388
$ws.Range("I93").Value = 800
$ws.Range("J93").Value = 55556624
$ws.Range("K93").Value = 800
$ws.Range("L93").Value = 55556624
$ws.Range("M93").Value = 448
$ws.Range("N93").Value = -55559120
$ws.Range("H103").Value = 1000
$ws.Range("J103").Value = 1000
$ws.Range("L103").Value = 1000
$ws.Range("N103").Value = -3344
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1651.5714
$ws.Range("J96").Value = 1918.2
$ws.Range("L96").Value = 1918.2
$ws.Range("N96").Value = -4664.2
$ws.Range("H140").Value = 62300
$ws.Range("J140").Value = 62300
$ws.Range("L140").Value = 62300
$ws.Range("N140").Value = -72660
$ws.Range("H141").Value = 85428.75
$ws.Range("J141").Value = 85428.75
$ws.Range("L141").Value = 85428.75
$ws.Range("N141").Value = -95788.75
